$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.841.71'
$ws.Range("E2").Value = '  -0.92%  '

$ws.Range("D3").Value = '3.480.77'
$ws.Range("E3").Value = '  -0.65%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '604.07'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.30%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '144.25'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -3.26%  '

$ws.Range("D7").Value = '3.479.15'
$ws.Range("E7").Value = '  -0.67%  '

$ws.Range("E8").Value = '  -0.02%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.477'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.78%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.140'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -2.38%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '7.91'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +4.31%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.417'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.73%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000212'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -1.55%  '

$ws.Range("D14").Value = '4.069.73'
$ws.Range("E14").Value = '  -0.59%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '30.93'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -3.29%  '

$ws.Range("D16").Value = '3.478.90'
$ws.Range("E16").Value = '  -0.99%  '

$ws.Range("D17").Value = '66.924.19'
$ws.Range("E17").Value = '  -0.84%  '

$ws.Range("E18").Value = '  -0.27%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '10.59'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +6.30%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.25'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.64%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '15.27'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.92%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '428.65'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -3.89%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.600'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -3.90%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '79.32'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.33%  '

$ws.Range("E25").Value = '  +0.07%  '

$ws.Range("D26").Value = '3.620.17'
$ws.Range("E26").Value = '  -0.74%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.0000115'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -5.26%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '9.68'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.36%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '8.07'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -5.15%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.49'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.59%  '

$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.04%  '

$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.54'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -6.20%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.165'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -1.63%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '25.30'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.26%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.76'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -3.67%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '7.90'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -1.01%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '5.68'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -8.54%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.17%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '174.70'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.30%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.0887'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.13%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.30'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.98%  '

$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.891'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.92%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.97'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -12.90%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '46.26'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.01%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '27.48'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -11.25%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.20'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -6.94%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.28'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.99%  '

$ws.Range("B49").Value = 'SuiNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.977'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.12%  '

$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.36'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -4.70%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.243'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -2.65%  '
